# ARKCORR-30: Show Suspend queue
# Adds a new "Suspend queue" rule row (row 24) to the OnEnterQueue rules
# table on Sheet1, mirroring the existing queue rows (Intake, Fulfill,
# Supervisor Approval, Executive Approval, Release).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing rule row (23) down onto the new
# row (24) so the new row matches the style of its siblings (style id 16).
$ws.Range("B23:D23").Copy()
$ws.Range("B24:D24").PasteSpecial(-4122)

# Populate the new "Suspend queue" rule.
$ws.Cells.Item(24, 2).Value = "Suspend queue"
$ws.Cells.Item(24, 3).Value = "Suspend"
$ws.Cells.Item(24, 4).Value = "correspondence-extension-suspend-process"

# Update the view so the new row is the active selection, matching the
# author's final cursor position on save.
$ws.Activate()
$excel.Goto($ws.Range("B16"), $true) | Out-Null
$ws.Range("D24").Select() | Out-Null
